$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-1.034126353949444, 1.1118527362954551, 0.5690569587064784, 2.52671880263419)"
$ws.Range("C2").Value = "NIG(1.3287203721622989, 1.0107625990493978, 4.475671666169665, 5.867890327889075)"
$ws.Range("D2").Value = "NCT(3.0266367271342656, 1.3206370410235402, -0.007913798292498282, 2.620452703406555)"
$ws.Range("E2").Value = "JSU(-1.3371489800552223, 1.195000257472815, 1.9353175993145326, 4.917654882011699)"
